$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing formatting like trailing zeros, e.g. "42.10" -> 42.1) are forced to
# Text format first so the literal string is preserved, matching the source data.

$ws.Range("D2").Value = "93.984.44"
$ws.Range("E2").Value = "  -3.12%  "

$ws.Range("D3").Value = "3.435.99"
$ws.Range("E3").Value = "  +2.65%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.11"
$ws.Range("E5").Value = "  -5.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.62"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.395"
$ws.Range("E8").Value = "  -7.34%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.967"
$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("D11").Value = "3.439.54"
$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.10"
$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.197"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "93.791.71"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").Value = "4.086.10"
$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.35"
$ws.Range("E18").Value = "  -5.98%  "

$ws.Range("D19").Value = "3.438.67"
$ws.Range("E19").Value = "  +2.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.56"
$ws.Range("E20").Value = "  -0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.29"
$ws.Range("E21").Value = "  +4.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.499"
$ws.Range("E22").Value = "  -12.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "496.06"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.13"
$ws.Range("E24").Value = "  -6.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.56"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("E26").Value = "  -5.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.82"
$ws.Range("E27").Value = "  -6.42%  "

$ws.Range("D28").Value = "3.621.29"
$ws.Range("E28").Value = "  +2.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.90"
$ws.Range("E29").Value = "  -2.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.68"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("E32").Value = "  +7.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  -7.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.181"
$ws.Range("E34").Value = "  -4.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.18"
$ws.Range("E36").Value = "  +5.74%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "539.85"
$ws.Range("E38").Value = "  +5.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.61"
$ws.Range("E39").Value = "  -3.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  -4.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.926"
$ws.Range("E42").Value = "  +9.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.151"
$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.02"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0410"
$ws.Range("E46").Value = "  -6.39%  "

$ws.Range("E47").Value = "  -3.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.49"
$ws.Range("E48").Value = "  -4.42%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.14"
$ws.Range("E49").Value = "  +6.83%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.51"
$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.16"
$ws.Range("E51").Value = "  +1.17%  "
